# This sheet lists daily price observations for "Bruselas (repollito)" at
# "Femacal de La Calera". The commit ("Fruta / hortaliza, semanal") re-pulls
# the weekly snapshot, which re-orders the existing daily rows (2-35) — the
# same 34 observations reappear, just attached to different report dates
# (column D) and, since volume/price depend on the date, different
# Volumen/Precio values (columns J, K, L, M, P). Columns A, B, C, E-I, N, O,
# Q, R are identical on every row, so the net effect is a pure permutation
# of rows 2-35.
#
# new row R receives exactly the old contents (D, J, K, L, M, P) of old row
# $rowMap[R].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 35

# new_row -> old_row (1-based worksheet row numbers)
$rowMap = @{
    2 = 6;  3 = 22; 4 = 17; 5 = 2;  6 = 26; 7 = 5;  8 = 34; 9 = 29; 10 = 9;
    11 = 18; 12 = 27; 13 = 21; 14 = 28; 15 = 35; 16 = 23; 17 = 16; 18 = 32;
    19 = 20; 20 = 3;  21 = 19; 22 = 25; 23 = 14; 24 = 11; 25 = 24; 26 = 30;
    27 = 8;  28 = 4;  29 = 31; 30 = 10; 31 = 33; 32 = 7;  33 = 15; 34 = 13;
    35 = 12
}

# Snapshot the columns that vary per-row BEFORE any writes, since the
# permutation has to be applied simultaneously (several rows are part of
# the same reshuffle cycle).
$colLetters = @("D", "J", "K", "L", "M", "P")
$colIndex = @{ "D" = 4; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "P" = 16 }

$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($col in $colLetters) {
        $rowData[$col] = $ws.Cells.Item($r, $colIndex[$col]).Value2
    }
    $snapshot[$r] = $rowData
}

# Write each new row from the snapshot of its mapped source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $colLetters) {
        $ws.Cells.Item($r, $colIndex[$col]).Value = $srcData[$col]
    }
}
